# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# The table's run properties (b/i/strike/sz/u/color) keep the same
# logical values, but need to be re-emitted by the (now newer) OOXML
# writer. Re-apply every run's character formatting explicitly so the
# document is rewritten through the current engine.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row, Col, Bold, Italic, Underline, StrikeThrough, SizePt, Color (hex RRGGBB or $null)
$cells = @(
    @{Row=1; Col=2; B=$true;  I=$false; U=$false; S=$false; Sz=10; Color=$null}
    @{Row=1; Col=3; B=$false; I=$true;  U=$false; S=$false; Sz=10; Color=$null}
    @{Row=1; Col=4; B=$false; I=$false; U=$true;  S=$false; Sz=10; Color=$null}
    @{Row=1; Col=5; B=$false; I=$false; U=$false; S=$true;  Sz=10; Color=$null}

    @{Row=2; Col=1; B=$true;  I=$false; U=$false; S=$false; Sz=10; Color=$null}
    @{Row=2; Col=2; B=$true;  I=$false; U=$false; S=$false; Sz=5;  Color="ff007f"}
    @{Row=2; Col=3; B=$true;  I=$true;  U=$false; S=$false; Sz=6;  Color="007fff"}
    @{Row=2; Col=4; B=$true;  I=$false; U=$true;  S=$false; Sz=7;  Color="7fff00"}
    @{Row=2; Col=5; B=$true;  I=$false; U=$false; S=$true;  Sz=8;  Color="ff007f"}

    @{Row=3; Col=1; B=$false; I=$true;  U=$false; S=$false; Sz=10; Color=$null}
    @{Row=3; Col=2; B=$true;  I=$true;  U=$false; S=$false; Sz=9;  Color="007fff"}
    @{Row=3; Col=3; B=$false; I=$true;  U=$false; S=$false; Sz=10; Color="7fff00"}
    @{Row=3; Col=4; B=$false; I=$true;  U=$true;  S=$false; Sz=11; Color="ff007f"}
    @{Row=3; Col=5; B=$false; I=$true;  U=$false; S=$true;  Sz=12; Color="007fff"}

    @{Row=4; Col=1; B=$false; I=$false; U=$true;  S=$false; Sz=10; Color=$null}
    @{Row=4; Col=2; B=$true;  I=$false; U=$true;  S=$false; Sz=13; Color="7fff00"}
    @{Row=4; Col=3; B=$false; I=$true;  U=$true;  S=$false; Sz=14; Color="ff007f"}
    @{Row=4; Col=4; B=$false; I=$false; U=$true;  S=$false; Sz=15; Color="007fff"}
    @{Row=4; Col=5; B=$false; I=$false; U=$true;  S=$true;  Sz=16; Color="7fff00"}

    @{Row=5; Col=1; B=$false; I=$false; U=$false; S=$true;  Sz=10; Color=$null}
    @{Row=5; Col=2; B=$true;  I=$false; U=$false; S=$true;  Sz=17; Color="ff007f"}
    @{Row=5; Col=3; B=$false; I=$true;  U=$false; S=$true;  Sz=18; Color="007fff"}
    @{Row=5; Col=4; B=$false; I=$false; U=$true;  S=$true;  Sz=19; Color="7fff00"}
    @{Row=5; Col=5; B=$false; I=$false; U=$false; S=$true;  Sz=20; Color="ff007f"}
)

foreach ($item in $cells) {
    $cell = $t.Rows.Item($item.Row).Cells.Item($item.Col)
    $font = $cell.Range.Font

    $font.Size = $item.Sz
    $font.Bold = $item.B
    $font.Italic = $item.I
    if ($item.U) {
        $font.Underline = 1
    } else {
        $font.Underline = 0
    }
    $font.StrikeThrough = $item.S
    if ($item.Color) {
        $r = [Convert]::ToInt32($item.Color.Substring(0,2), 16)
        $g = [Convert]::ToInt32($item.Color.Substring(2,2), 16)
        $b = [Convert]::ToInt32($item.Color.Substring(4,2), 16)
        $font.Color = $r + ($g * 256) + ($b * 65536)
    }
}
